$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ExtremeFlows")

# Fix typo: missing space after "C." in A4
$ws.Range("A4").Value = "C. Low Lake Powell releases + gains through Grand Canyon"

# Replace leading-space labels with properly prefixed "F1./F2./F3." labels
$ws.Range("A8").Value = "F1. 10-year"
$ws.Range("A9").Value = "F2. 4-year"
$ws.Range("A10").Value = "F3. 3-year"

# Update the selected cell on the sheet so printing/viewing starts at A7
$ws.Range("A7").Select()
